$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Run the whole suite: flip Runmode (column C) to "Y" for every test case
# except TestCase_A13 (row 14), which flips to "N".
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("C10").Value = "Y"
$ws.Range("C11").Value = "Y"
$ws.Range("C12").Value = "Y"
$ws.Range("C13").Value = "Y"
$ws.Range("C14").Value = "N"
$ws.Range("C15").Value = "Y"
$ws.Range("C16").Value = "Y"
$ws.Range("C17").Value = "Y"
$ws.Range("C18").Value = "Y"

# Leave the cursor where the author left it when they saved.
$ws.Range("C15").Select()
